$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.144.47"
$ws.Range("E2").Value = "  +2.41%  "

# Row 3
$ws.Range("D3").Value = "2.983.74"
$ws.Range("E3").Value = "  +1.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.21"
$ws.Range("E5").Value = "  +1.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.41"
$ws.Range("E6").Value = "  +7.52%  "

# Row 7
$ws.Range("E7").Value = "  -0.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  +2.19%  "

# Row 9
$ws.Range("D9").Value = "2.973.48"
$ws.Range("E9").Value = "  +0.98%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  +4.02%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.15"
$ws.Range("E11").Value = "  +6.77%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  +2.00%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +4.74%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.70"
$ws.Range("E14").Value = "  +2.96%  "

# Row 15
$ws.Range("E15").Value = "  +2.37%  "

# Row 16
$ws.Range("D16").Value = "3.478.61"
$ws.Range("E16").Value = "  +1.37%  "

# Row 17
$ws.Range("E17").Value = "  +6.99%  "

# Row 18
$ws.Range("D18").Value = "2.983.14"
$ws.Range("E18").Value = "  +0.94%  "

# Row 19
$ws.Range("D19").Value = "59.124.19"
$ws.Range("E19").Value = "  +2.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "427.92"
$ws.Range("E20").Value = "  +3.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("E21").Value = "  +4.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  +4.90%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("E23").Value = "  +2.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.41"
$ws.Range("E24").Value = "  +4.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.69"
$ws.Range("E25").Value = "  +1.91%  "

# Row 26
$ws.Range("E26").Value = "  -0.13%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.16%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.15"
$ws.Range("E28").Value = "  +9.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.54"
$ws.Range("E29").Value = "  +2.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.70"
$ws.Range("E30").Value = "  +2.84%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.69"
$ws.Range("E31").Value = "  +2.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.09"
$ws.Range("E32").Value = "  -0.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0986"
$ws.Range("E33").Value = "  -4.11%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0771"
$ws.Range("E34").Value = "  +17.89%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.991"
$ws.Range("E35").Value = "  +6.15%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.85"
$ws.Range("E36").Value = "  +4.51%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.08"
$ws.Range("E37").Value = "  +0.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.27"
$ws.Range("E38").Value = "  +1.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.63"
$ws.Range("E39").Value = "  +3.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("E40").Value = "  +7.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "400.37"
$ws.Range("E41").Value = "  +7.14%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.766.02"
$ws.Range("E42").Value = "  +4.54%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0351"
$ws.Range("E43").Value = "  +1.87%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.108"
$ws.Range("E44").Value = "  +0.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.250"
$ws.Range("E45").Value = "  +6.35%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.23"
$ws.Range("E47").Value = "  +0.56%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.110"
$ws.Range("E48").Value = "  +1.23%  "

# Row 49
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.92"
$ws.Range("E49").Value = "  +19.02%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("E50").Value = "  +0.93%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.40"
$ws.Range("E51").Value = "  +0.78%  "
